$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C4").Value = -12.901
$ws.Range("A9").Value = -21.76
$ws.Range("C9").Value = -12.42
$ws.Range("D9").Value = -7.936999999999999
$ws.Range("C11").Value = -12.14
$ws.Range("A13").Value = -21.884
$ws.Range("A16").Value = -21.868
$ws.Range("C16").Value = -13.073
$ws.Range("A18").Value = -21.867
$ws.Range("A20").Value = -20.762
$ws.Range("D22").Value = -7.805
$ws.Range("C23").Value = -12.927
$ws.Range("C24").Value = -12.646
$ws.Range("A26").Value = -21.839
$ws.Range("C26").Value = -12.904
$ws.Range("A27").Value = -21.78
$ws.Range("D27").Value = -7.762
$ws.Range("A29").Value = -21.153
$ws.Range("D29").Value = -7.706
$ws.Range("D32").Value = -7.609
$ws.Range("C34").Value = -12.579
$ws.Range("A35").Value = -20.339
$ws.Range("C35").Value = -12.58
$ws.Range("A36").Value = -20.266
$ws.Range("D37").Value = -7.715000000000001
$ws.Range("D38").Value = -7.650999999999999
$ws.Range("D39").Value = -7.479000000000001
$ws.Range("D41").Value = -8.215
$ws.Range("C44").Value = -12.282
$ws.Range("A45").Value = -21.519
$ws.Range("D45").Value = -7.528
$ws.Range("C48").Value = -12.986
$ws.Range("D48").Value = -7.787999999999999
$ws.Range("C49").Value = -12.907
$ws.Range("D51").Value = -7.970000000000001
$ws.Range("C52").Value = -12.095
$ws.Range("A55").Value = -22.102
$ws.Range("D56").Value = -8.381
$ws.Range("A57").Value = -21.962
$ws.Range("D57").Value = -8.016
$ws.Range("D61").Value = -7.433
$ws.Range("D64").Value = -7.931999999999999
$ws.Range("C66").Value = -12.076
$ws.Range("C67").Value = -11.554
$ws.Range("A69").Value = -21.466
$ws.Range("C73").Value = -12.088
$ws.Range("D75").Value = -7.882000000000001
$ws.Range("A76").Value = -20.347
$ws.Range("A78").Value = -20.045
$ws.Range("C78").Value = -12.378
$ws.Range("C80").Value = -12.324
$ws.Range("A82").Value = -21.801
$ws.Range("D82").Value = -8.023
$ws.Range("A83").Value = -21.84
$ws.Range("D90").Value = -7.263000000000001
$ws.Range("C91").Value = -12.856
$ws.Range("A93").Value = -21.428
$ws.Range("D93").Value = -7.441
$ws.Range("A97").Value = -21.719
$ws.Range("C97").Value = -11.58
$ws.Range("C99").Value = -12.217
$ws.Range("D102").Value = -7.656000000000001
$ws.Range("C104").Value = -12.906
$ws.Range("D105").Value = -7.721000000000001